$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit2"
$ws.Range("C2").Value = "Sdc1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.03343333333333334
$ws.Range("H2").Value = 0.1003
$ws.Range("I2").Value = 0.01753259568243662
$ws.Range("J2").Value = 0.01753259568243662
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8369776666666665
$ws.Range("N2").Value = 2.510933
$ws.Range("O2").Value = 0.0694586718035551
$ws.Range("P2").Value = 0.06945867180355511
$ws.Range("Q2").Value = 0.02798295332222222
$ws.Range("R2").Value = 0.2518465799
$ws.Range("S2").Value = 0.001217790809370792
$ws.Range("T2").Value = 0.001217790809370792
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit2"
$ws.Range("C3").Value = "Sdc1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.03343333333333334
$ws.Range("H3").Value = 0.1003
$ws.Range("I3").Value = 0.01753259568243662
$ws.Range("J3").Value = 0.01753259568243662
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.427350333333333
$ws.Range("N3").Value = 7.282051
$ws.Range("O3").Value = 0.2014397000898671
$ws.Range("P3").Value = 0.2014397000898671
$ws.Range("Q3").Value = 0.08115441281111112
$ws.Range("R3").Value = 0.7303897153
$ws.Range("S3").Value = 0.003531760816066933
$ws.Range("T3").Value = 0.003531760816066932
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit2"
$ws.Range("C4").Value = "Sdc1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.03343333333333334
$ws.Range("H4").Value = 0.1003
$ws.Range("I4").Value = 0.01753259568243662
$ws.Range("J4").Value = 0.01753259568243662
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.785681666666667
$ws.Range("N4").Value = 26.357045
$ws.Range("O4").Value = 0.7291016281065776
$ws.Range("P4").Value = 0.7291016281065776
$ws.Range("Q4").Value = 0.2937346237222223
$ws.Range("R4").Value = 2.6436116135
$ws.Range("S4").Value = 0.0127830440569989
$ws.Range("T4").Value = 0.01278304405699889
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slit2"
$ws.Range("C5").Value = "Sdc1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.37963
$ws.Range("H5").Value = 1.13889
$ws.Range("I5").Value = 0.1990797397484571
$ws.Range("J5").Value = 0.1990797397484571
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8369776666666665
$ws.Range("N5").Value = 2.510933
$ws.Range("O5").Value = 0.0694586718035551
$ws.Range("P5").Value = 0.06945867180355511
$ws.Range("Q5").Value = 0.3177418315966666
$ws.Range("R5").Value = 2.85967648437
$ws.Range("S5").Value = 0.01382781430592524
$ws.Range("T5").Value = 0.01382781430592524
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit2"
$ws.Range("C6").Value = "Sdc1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.37963
$ws.Range("H6").Value = 1.13889
$ws.Range("I6").Value = 0.1990797397484571
$ws.Range("J6").Value = 0.1990797397484571
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.427350333333333
$ws.Range("N6").Value = 7.282051
$ws.Range("O6").Value = 0.2014397000898671
$ws.Range("P6").Value = 0.2014397000898671
$ws.Range("Q6").Value = 0.9214950070433333
$ws.Range("R6").Value = 8.293455063390001
$ws.Range("S6").Value = 0.04010256306889799
$ws.Range("T6").Value = 0.04010256306889799
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit2"
$ws.Range("C7").Value = "Sdc1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.37963
$ws.Range("H7").Value = 1.13889
$ws.Range("I7").Value = 0.1990797397484571
$ws.Range("J7").Value = 0.1990797397484571
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.785681666666667
$ws.Range("N7").Value = 26.357045
$ws.Range("O7").Value = 0.7291016281065776
$ws.Range("P7").Value = 0.7291016281065776
$ws.Range("Q7").Value = 3.335308331116666
$ws.Range("R7").Value = 30.01777498005
$ws.Range("S7").Value = 0.1451493623736338
$ws.Range("T7").Value = 0.1451493623736338
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Slit2"
$ws.Range("C8").Value = "Sdc1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.493861
$ws.Range("H8").Value = 4.481583000000001
$ws.Range("I8").Value = 0.7833876645691064
$ws.Range("J8").Value = 0.7833876645691064
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8369776666666665
$ws.Range("N8").Value = 2.510933
$ws.Range("O8").Value = 0.0694586718035551
$ws.Range("P8").Value = 0.06945867180355511
$ws.Range("Q8").Value = 1.250328294104333
$ws.Range("R8").Value = 11.252954646939
$ws.Range("S8").Value = 0.05441306668825906
$ws.Range("T8").Value = 0.05441306668825908
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Slit2"
$ws.Range("C9").Value = "Sdc1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.493861
$ws.Range("H9").Value = 4.481583000000001
$ws.Range("I9").Value = 0.7833876645691064
$ws.Range("J9").Value = 0.7833876645691064
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.427350333333333
$ws.Range("N9").Value = 7.282051
$ws.Range("O9").Value = 0.2014397000898671
$ws.Range("P9").Value = 0.2014397000898671
$ws.Range("Q9").Value = 3.626123996303667
$ws.Range("R9").Value = 32.635115966733
$ws.Range("S9").Value = 0.1578053762049022
$ws.Range("T9").Value = 0.1578053762049022
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Slit2"
$ws.Range("C10").Value = "Sdc1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.493861
$ws.Range("H10").Value = 4.481583000000001
$ws.Range("I10").Value = 0.7833876645691064
$ws.Range("J10").Value = 0.7833876645691064
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.785681666666667
$ws.Range("N10").Value = 26.357045
$ws.Range("O10").Value = 0.7291016281065776
$ws.Range("P10").Value = 0.7291016281065776
$ws.Range("Q10").Value = 13.12458720024834
$ws.Range("R10").Value = 118.121284802235
$ws.Range("S10").Value = 0.571169221675945
$ws.Range("T10").Value = 0.571169221675945
